# Apply "Final commit with lot of minor errors solved" edits:
#  - Fix misspelled name "Prretika Shetty" -> "Preetika Shetty" in row 3 (col B)
#  - Column A: change id numbers 1..10 -> 52501..52510
#  - Column C: update several score values
#  - Selection changes from whole range A1:C10 to single active cell B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ID values in column A
$idValues = @(52501,52502,52503,52504,52505,52506,52507,52508,52509,52510)

# New names in column B (row3 corrected spelling, others unchanged)
$names = @(
    "Jagannath Pidaparthy",
    "Vishal Patil",
    "Preetika Shetty",
    "Sagar Mishra",
    "Shubham Mishra",
    "Kanchan Soni",
    "Jai Lohani",
    "Korol Dhanda",
    "kaustubh Srivastava",
    "Purva Shinde"
)

# New score values in column C
$scores = @(80,80,97,96,93,95,92,85,91,90)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $idValues[$i]
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $scores[$i]
}

# Update the selection: single active cell B3 instead of the A1:C10 range
$ws.Range("B3").Select()
